$wb = $excel.ActiveWorkbook

# Deckel das Temperaturfeld (Anzahl der Oszillationen / C2) auf die
# Verdampfungstemperatur-bedingte Obergrenze -> Wert von 100 auf 5 senken.
$ws5 = $wb.Worksheets.Item("Tabelle5")
$ws5.Range("C2").Value = 5

# Tabelle5 wird zum aktiven Blatt; dadurch wandert auch die Selektion
# (tabSelected) von Tabelle2 auf Tabelle5, und die Auswahl zeigt auf C2.
$ws5.Activate()
$ws5.Range("C2").Select()
